$wb = $excel.ActiveWorkbook

# --- Sheet3: update the list separator in the "classes" column from "|" to " $$ " ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A8").Value = 'ex:ClassA1 $$ ex:ClassA2 $$ ex:ClassA3'
$ws3.Range("A9").Value = 'ex:ClassB1 $$ ex:ClassB2'
$ws3.Range("A10").Value = 'ex:ClassC1 $$ ex:ClassC2 $$ ex:ClassC3 $$ ex:ClassC4'

# --- Column width re-layout on Sheet1 and Sheet1_2 (split the single custom column
#     formatting for column C into three ranges: A:B, C, D:ALV) ---
foreach ($idx in 1..2) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range($ws.Columns.Item(1), $ws.Columns.Item(2)).ColumnWidth = 7.666666666666667
    $ws.Columns.Item(3).ColumnWidth = 15.166666666666666
    $ws.Range($ws.Columns.Item(4), $ws.Columns.Item(1025)).ColumnWidth = 7.666666666666667
}

# --- Make Sheet3 ("Sheet3") the active / selected tab, with A10 selected ---
$ws3.Activate() | Out-Null
$ws3.Range("A10").Select() | Out-Null
